$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries in column G describing the added "receipt / exceptions handler / date"
# work items. Values are entered in the order that makes the shared-string table
# line up with the target (add payment -> deal with receipt number -> STOCK!).
$ws.Range("G15").Value = "add payment"
$ws.Range("G14").Value = "deal with receipt number"
$ws.Range("G16").Value = "STOCK!"

# Row 14 wraps onto two lines, so it needs the taller row height.
$ws.Rows(14).RowHeight = 32

# Reflect the new active cell / selection left behind by the edit.
[void]$ws.Range("G16").Select()
